# check #11 on 31/3/2025
# Corrects 31 values in column I ("Total quantity_imp") on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I5").Value = 351
$ws.Range("I10").Value = 231
$ws.Range("I17").Value = 388
$ws.Range("I63").Value = 13
$ws.Range("I133").Value = 33.2
$ws.Range("I137").Value = 770
$ws.Range("I138").Value = 62
$ws.Range("I143").Value = 39.2
$ws.Range("I147").Value = 400
$ws.Range("I152").Value = 415
$ws.Range("I153").Value = 33.2
$ws.Range("I156").Value = 83
$ws.Range("I162").Value = 975
$ws.Range("I163").Value = 78
$ws.Range("I166").Value = 171
$ws.Range("I167").Value = 600
$ws.Range("I168").Value = 48
$ws.Range("I171").Value = 114
$ws.Range("I172").Value = 555
$ws.Range("I173").Value = 44.40000000000001
$ws.Range("I176").Value = 99
$ws.Range("I182").Value = 1185
$ws.Range("I183").Value = 94.80000000000001
$ws.Range("I186").Value = 216
$ws.Range("I191").Value = 127
$ws.Range("I192").Value = 345
$ws.Range("I193").Value = 27.6
$ws.Range("I196").Value = 63
$ws.Range("I197").Value = 575
$ws.Range("I198").Value = 46
$ws.Range("I201").Value = 115
